## Update "Forecast Comparison" sheet with the corrected forecast output:
##  - insert a new "Week_Start_Date" column after "Week"
##  - shorten the week labels (W01 -> W1 ... W09 -> W9)
##  - correct the MyForecast figures
##  - store is_holiday_week as a real boolean
## and refresh the dependent totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")
$summary = $wb.Worksheets.Item("Summary")

# --- Insert the new "Week_Start_Date" column as column B -----------------
$ws.Columns.Item(2).Insert()

$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# --- Per-week data: Week label, Week start date, corrected MyForecast ----
$weeks = @(
    @{ Row = 2;  Label = "W1";  Date = "2025-01-05"; Forecast = 70 },
    @{ Row = 3;  Label = "W2";  Date = "2025-01-12"; Forecast = 79 },
    @{ Row = 4;  Label = "W3";  Date = "2025-01-19"; Forecast = 77 },
    @{ Row = 5;  Label = "W4";  Date = "2025-01-26"; Forecast = 71 },
    @{ Row = 6;  Label = "W5";  Date = "2025-02-02"; Forecast = 76 },
    @{ Row = 7;  Label = "W6";  Date = "2025-02-09"; Forecast = 81 },
    @{ Row = 8;  Label = "W7";  Date = "2025-02-16"; Forecast = 80 },
    @{ Row = 9;  Label = "W8";  Date = "2025-02-23"; Forecast = 75 },
    @{ Row = 10; Label = "W9";  Date = "2025-03-02"; Forecast = 78 },
    @{ Row = 11; Label = "W10"; Date = "2025-03-09"; Forecast = 86 },
    @{ Row = 12; Label = "W11"; Date = "2025-03-16"; Forecast = 85 },
    @{ Row = 13; Label = "W12"; Date = "2025-03-23"; Forecast = 77 },
    @{ Row = 14; Label = "W13"; Date = "2025-03-30"; Forecast = 80 },
    @{ Row = 15; Label = "W14"; Date = "2025-04-06"; Forecast = 86 },
    @{ Row = 16; Label = "W15"; Date = "2025-04-13"; Forecast = 85 },
    @{ Row = 17; Label = "W16"; Date = "2025-04-20"; Forecast = 79 }
)

foreach ($week in $weeks) {
    $r = $week.Row

    # Week label (A) - force text so e.g. "W1" never becomes numeric
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $week.Label

    # Week_Start_Date (B) - keep as plain text, not an auto-converted date serial
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $week.Date

    # MyForecast (D, after the column insert shifted C -> D)
    $ws.Cells.Item($r, 4).Value = $week.Forecast

    # is_holiday_week (J, after the column insert shifted I -> J) - store as boolean
    $ws.Cells.Item($r, 10).Value = $false
}

# --- Refresh the dependent Summary sheet totals ---------------------------
# These are stored as text (matching the rest of the "Value" column), so
# force text formatting before assigning to avoid Excel coercing them to
# numbers.
$summaryCells = @(9, 10, 11, 12, 14)
foreach ($row in $summaryCells) {
    $summary.Cells.Item($row, 2).NumberFormat = "@"
}

$summary.Range("B9").Value  = "1266"   # Total Forecast (16 Weeks)
$summary.Range("B10").Value = "609"    # Total Forecast (8 Weeks)
$summary.Range("B11").Value = "298"    # Total Forecast (4 Weeks)
$summary.Range("B12").Value = "86"     # Max Forecast
$summary.Range("B14").Value = "70"     # Min Forecast
